# Ajustes A la Nomina
# - Rename sheet from "Sheet1" to "Datos"
# - Remove the "Revisar" column (H) header/formatting entirely
# - Strip the bold + centered header-row styling back to normal

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet.
$ws.Name = "Datos"

# Drop the now-unused "Revisar" column (H) completely.
$ws.Columns.Item(8).Delete()

# Remove the bold / center-aligned styling that was applied to the header row,
# returning it to the workbook's default (unstyled) formatting.
$ws.Rows.Item(1).ClearFormats()
